$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet positioned right after Sheet1, and name it "Result"
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Result"

# Populate the new "Result" sheet with the lookup table used to
# replace the removed CheckInputsAtOnce/etc. helper methods.
$ws.Range("A1").Value = "inputName"
$ws.Range("B1").Value = "isPassed"

$ws.Range("A2").Value = "IRC nick"
$ws.Range("B2").Value = "passed"

$ws.Range("A3").Value = "Email *"
$ws.Range("B3").Value = "failed"

$ws.Range("A4").Value = "Имя *"
$ws.Range("B4").Value = "passed"

# Keep Sheet1 as the active/selected tab, as before the edit.
$ws1.Activate()
